$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = 0.01339202195363811
$ws.Range("C2").Value = 0.751747495968654
$ws.Range("D2").Value = 1.207807154213344
$ws.Range("E2").Value = 1.099002799911512
$ws.Range("F2").Value = 1.111338784757336

# Row 3 (Q0)
$ws.Range("B3").Value = 0.1058108245021329
$ws.Range("C3").Value = 1.129124555725148
$ws.Range("D3").Value = 3.360206299466695
$ws.Range("E3").Value = 1.833086549911568
$ws.Range("F3").Value = 1.836417714414163
$ws.Range("G3").Value = 144

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1993574677006487
$ws.Range("C4").Value = 1.286806921283597
$ws.Range("D4").Value = 7.446533758047825
$ws.Range("E4").Value = 2.728833772520383
$ws.Range("F4").Value = 2.74119228546971
$ws.Range("G4").Value = 70
